$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.983317657231666, 50.06981459567332]"
$ws.Range("T2").Value = "[49.98509068464103, 50.04632671179001]"
$ws.Range("L3").Value = "[49.90765865717907, 50.088319966306926]"
$ws.Range("T3").Value = "[49.95826818553987, 50.05334269633412]"
